$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain stored as plain text,
# matching the inline-string format used in the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.103.46"
$ws.Range("E2").Value = "  -4.43%  "

$ws.Range("D3").Value = "2.975.50"
$ws.Range("E3").Value = "  -1.51%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "558.62"
$ws.Range("E5").Value = "  -3.62%  "

$ws.Range("D6").Value = "133.78"
$ws.Range("E6").Value = "  +4.55%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +2.97%  "

$ws.Range("D9").Value = "2.975.28"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("E10").Value = "  -3.32%  "

$ws.Range("E11").Value = "  -6.20%  "

$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").Value = "33.06"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("D16").Value = "3.465.83"
$ws.Range("E16").Value = "  -1.29%  "

$ws.Range("D17").Value = "6.91"
$ws.Range("E17").Value = "  +7.06%  "

$ws.Range("D18").Value = "2.973.35"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "57.999.00"
$ws.Range("E19").Value = "  -4.15%  "

$ws.Range("D20").Value = "420.99"
$ws.Range("E20").Value = "  -3.28%  "

$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("D23").Value = "7.02"
$ws.Range("E23").Value = "  -0.87%  "

$ws.Range("D24").Value = "13.13"
$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("D25").Value = "79.70"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").Value = "2.51"
$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("D29").Value = "7.61"
$ws.Range("E29").Value = "  +2.98%  "

$ws.Range("E30").Value = "  +4.71%  "

$ws.Range("D31").Value = "25.34"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("E32").Value = "  -2.24%  "

$ws.Range("E33").Value = "  +5.55%  "

$ws.Range("D34").Value = "2.15"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").Value = "0.945"
$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("D37").Value = "0.0₃0701"
$ws.Range("E37").Value = "  +3.49%  "

$ws.Range("E38").Value = "  -2.90%  "

$ws.Range("D39").Value = "8.67"
$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").Value = "2.59"
$ws.Range("E40").Value = "  +3.15%  "

$ws.Range("D41").Value = "0.0352"
$ws.Range("E41").Value = "  -2.93%  "

$ws.Range("D42").Value = "380.60"
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("E43").Value = "  -2.76%  "

$ws.Range("D44").Value = "2.691.75"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D46").Value = "0.244"
$ws.Range("E46").Value = "  +2.29%  "

$ws.Range("D47").Value = "122.11"
$ws.Range("E47").Value = "  +3.11%  "

$ws.Range("E48").Value = "  +2.64%  "

$ws.Range("D49").Value = "2.00"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").Value = "23.64"
$ws.Range("E50").Value = "  -1.76%  "

$ws.Range("E51").Value = "  -0.95%  "
